$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 598, shifting existing rows 598:630 down to 599:631.
$ws.Rows(598).Insert()

# Populate the newly inserted row 598 with the new price observation.
$ws.Range("A598").Value = 8
$ws.Range("B598").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C598").Value = 'Coquimbo'
$ws.Range("D598").Value = 45267
$ws.Range("E598").Value = 4
$ws.Range("F598").Value = 100112032
$ws.Range("G598").Value = 'Zapallo italiano'
$ws.Range("H598").Value = 'Sin especificar'
$ws.Range("I598").Value = 'Primera'
$ws.Range("J598").Value = 500
$ws.Range("K598").Value = 8000
$ws.Range("L598").Value = 9000
$ws.Range("M598").Value = 8500
$ws.Range("N598").Value = '$/caja 60 unidades'
$ws.Range("O598").Value = 'Provincia de Limarí'
$ws.Range("P598").Value = 142
$ws.Range("Q598").Value = 60
$ws.Range("R598").Value = 'Hortaliza'
